# Schema update for table_schema sheet:
#  - Insert a new column D ("is_nullable"), shifting the old
#    table_description/column_description columns from D/E to E/F.
#  - Rename column_type -> data_type (column C header).
#  - Populate the new is_nullable column header with the same bold/
#    centered/top-aligned header style as the other headers, but without
#    the cell border (matches the new cellXfs entry in the target file).
#  - Fill is_nullable data values (YES/NO) for each of the 41 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table_schema")

# Insert a new column before D; this shifts the existing
# table_description (D) / column_description (E) columns one to the
# right (-> E / F) and carries the header formatting along.
$ws.Columns.Item(4).Insert()

# Rename the column_type header to data_type (style/position unchanged).
$ws.Range("C1").Value = "data_type"

# Set the new column header text.
$ws.Range("D1").Value = "is_nullable"

# The inserted column copied the bordered header style from column C.
# The target style for this header has no cell border, so clear it -
# this creates the new (3rd) cellXfs entry: bold font, center/top
# alignment, no border - matching fontId=1/borderId=0 in the diff.
$ws.Range("D1").Borders.LineStyle = -4142

# is_nullable values for rows 2..42 (one per schema column row).
$isNullable = @("NO","NO","YES","NO","NO","NO","NO","NO","NO","NO","NO","NO","YES","NO","NO","NO","NO","YES","NO","NO","NO","NO","NO","YES","YES","NO","NO","NO","NO","NO","NO","NO","YES","NO","NO","NO","NO","YES","YES","YES","YES")

for ($i = 0; $i -lt $isNullable.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $isNullable[$i]
}
